$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.733.08"
$ws.Range("E2").Value = "  +1.83%  "
$ws.Range("D3").Value = "1.574.58"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "213.83"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("D6").Value = "0.492"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "44.71"
$ws.Range("E8").Value = "  +2.10%  "
$ws.Range("D9").Value = "24.28"
$ws.Range("E9").Value = "  +1.60%  "
$ws.Range("E10").Value = "  -1.09%  "
$ws.Range("E11").Value = "  -1.08%  "
$ws.Range("D12").Value = "0.0889"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").Value = "1.800.13"
$ws.Range("E13").Value = "  -0.87%  "
$ws.Range("D14").Value = "1.574.97"
$ws.Range("E14").Value = "  -0.92%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "0.523"
$ws.Range("E15").Value = "  -1.34%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "28.723.82"
$ws.Range("E16").Value = "  +1.63%  "
$ws.Range("D17").Value = "3.68"
$ws.Range("E17").Value = "  -1.62%  "
$ws.Range("D18").Value = "62.54"
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("D19").Value = "230.59"
$ws.Range("E19").Value = "  +1.59%  "
$ws.Range("E20").Value = "  -0.83%  "
$ws.Range("D21").Value = "0.0₃0694"
$ws.Range("E21").Value = "  -1.90%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("E23").Value = "  -4.83%  "
$ws.Range("E24").Value = "  -1.39%  "
$ws.Range("D25").Value = "2.11"
$ws.Range("E25").Value = "  +8.30%  "
$ws.Range("D26").Value = "152.16"
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("E28").Value = "  -1.71%  "
$ws.Range("E29").Value = "  -2.45%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  +2.74%  "
$ws.Range("E32").Value = "  -1.93%  "
$ws.Range("E33").Value = "  -0.46%  "
$ws.Range("E34").Value = "  -1.21%  "
$ws.Range("D35").Value = "1.391.84"
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("E36").Value = "  +2.98%  "
$ws.Range("D37").Value = "1.55"
$ws.Range("E37").Value = "  -2.97%  "
$ws.Range("D38").Value = "2.37"
$ws.Range("E38").Value = "  +0.68%  "
$ws.Range("E39").Value = "  +2.83%  "
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("D41").Value = "0.526"
$ws.Range("E41").Value = "  -2.49%  "
$ws.Range("E42").Value = "  +1.94%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "0.795"
$ws.Range("E44").Value = "  -2.15%  "
$ws.Range("E45").Value = "  +3.11%  "
$ws.Range("D46").Value = "5.53"
$ws.Range("E46").Value = "  -1.08%  "
$ws.Range("D47").Value = "0.962"
$ws.Range("E47").Value = "  -1.96%  "
$ws.Range("D48").Value = "63.34"
$ws.Range("E48").Value = "  -1.36%  "
$ws.Range("D49").Value = "1.711.90"
$ws.Range("E49").Value = "  -0.71%  "
$ws.Range("D50").Value = "86.86"
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.0517"
$ws.Range("E51").Value = "  -1.07%  "
